$wb = $excel.ActiveWorkbook

# --- Rename "Sheet3" to "Device Review" and populate it with the new
#     RF amplifier device table (Device / Part Number / Vendor / Description) ---
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Name = "Device Review"

$ws.Range("C5").Value = "Device"
$ws.Range("D5").Value = "Part Number"
$ws.Range("E5").Value = "Vendor"
$ws.Range("F5").Value = "Description"

$ws.Range("C6").Value = "RF AMP"
$ws.Range("D6").Value = "MRFE6VP61K25H"
$ws.Range("E6").Value = "NXP"
$ws.Range("F6").Value = "1250W RF Power LDMOS Transistor"

# Widen the Part Number / Description columns to fit their contents
# (mirrors the workbook's original "best fit" column widths).
$ws.Columns("D").ColumnWidth = 16
$ws.Columns("F").ColumnWidth = 33.29

# Put the selection where the author left it and make this the active sheet/tab.
$ws.Range("P11").Select() | Out-Null
$ws.Activate() | Out-Null
